$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7445297241210938
$ws.Range("B1").Value = 1.582640290260315
$ws.Range("C1").Value = 4.748236179351807
$ws.Range("D1").Value = 2.402153491973877
$ws.Range("E1").Value = 1.269577145576477
